$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # model_characteristics
$ws2 = $wb.Worksheets.Item(2)   # model_overview

# --- Correction of Fujita and Swameye parameter counts ---

# model_characteristics sheet (Fujita = row 11, Swameye = row 19)
# Parameters column F: Fujita 22 -> 19, Swameye 16 -> 13
$ws1.Range("F11").Value = 19
$ws1.Range("F19").Value = 13
# Error pars column I: Fujita 3 -> 0, Swameye 3 -> 0
$ws1.Range("I11").Value = 0
$ws1.Range("I19").Value = 0

# model_overview sheet (Fujita = row 12, Swameye = row 20)
# Error pars column I: 3 -> 0
$ws2.Range("I12").Value = 0
$ws2.Range("I20").Value = 0
# Total estimated column J: now a plain corrected value instead of the G+H+I formula
$ws2.Range("J12").Value = 19
$ws2.Range("J20").Value = 13

# --- Update sheet selections / active sheet ---

# On model_overview, move the selection to H12
$ws2.Activate()
$ws2.Range("H12").Select()

# Make model_characteristics the active sheet with selection on G19
$ws1.Activate()
$ws1.Range("G19").Select()
